$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 424: the "C424" cell used to contain the text "NA"; the latest
# script run found no real page number for that day, so the value is
# cleared (kept as an empty/blank cell, same as neighbouring rows 422/423).
$ws.Range("C424").Value = ""

# New rows appended by the latest script run (2026-01-30 results).
$newRows = @(
    @(425, "2026-01-30", "bonnes pratiques", 31, 1),
    @(426, "2026-01-30", "espèces exotiques envahissantes", 83, 1),
    @(427, "2026-01-30", "substance active", 84, 2),
    @(428, "2026-01-30", "eaux souterraines", 84, 1),
    @(429, "2026-01-30", "eaux de surface", 84, 1),
    @(430, "2026-01-30", "substances actives", 84, 1),
    @(431, "2026-01-30", "insecticide", 85, 1),
    @(432, "2026-01-30", "substance active", 85, 1)
)

foreach ($row in $newRows) {
    $r = $row[0]
    # Column A holds a date-like string ("2026-01-30"); entering it as a
    # literal formula that evaluates to that text, then converting the
    # formula to a static value, keeps it stored as plain text (matching
    # every other row in the sheet) instead of Excel auto-converting it
    # into a date serial number/format.
    $ws.Cells.Item($r, 1).Formula = "=""" + $row[1] + """"
    $ws.Cells.Item($r, 2).Value = $row[2]
    $ws.Cells.Item($r, 3).Value = $row[3]
    $ws.Cells.Item($r, 4).Value = $row[4]
}

$ws.Range("A425:A432").Copy()
$ws.Range("A425:A432").PasteSpecial(-4163)
